$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.04575233333333333
$ws.Range("H2").Value = 0.137257
$ws.Range("I2").Value = 0.14147347546269
$ws.Range("J2").Value = 0.14147347546269
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.001065388834
$ws.Range("R2").Value = 0.009588499505999999
$ws.Range("S2").Value = 0.001317218570037659
$ws.Range("T2").Value = 0.001317218570037659

# Row 3
$ws.Range("G3").Value = 0.04575233333333333
$ws.Range("H3").Value = 0.137257
$ws.Range("I3").Value = 0.14147347546269
$ws.Range("J3").Value = 0.14147347546269
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.00608136964511111
$ws.Range("R3").Value = 0.05473232680599999
$ws.Range("S3").Value = 0.007518844549673292
$ws.Range("T3").Value = 0.007518844549673291

# Row 4
$ws.Range("G4").Value = 0.04575233333333333
$ws.Range("H4").Value = 0.137257
$ws.Range("I4").Value = 0.14147347546269
$ws.Range("J4").Value = 0.14147347546269
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.1072794001657778
$ws.Range("R4").Value = 0.9655146014920001
$ws.Range("S4").Value = 0.1326374123429791
$ws.Range("T4").Value = 0.132637412342979

# Row 5
$ws.Range("I5").Value = 0.2599907647526892
$ws.Range("J5").Value = 0.2599907647526892
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.001957902404
$ws.Range("R5").Value = 0.017621121636
$ws.Range("S5").Value = 0.002420698736992935
$ws.Range("T5").Value = 0.002420698736992936

# Row 6
$ws.Range("I6").Value = 0.2599907647526892
$ws.Range("J6").Value = 0.2599907647526892
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.01381764417770089
$ws.Range("T6").Value = 0.01381764417770089

# Row 7
$ws.Range("I7").Value = 0.2599907647526892
$ws.Range("J7").Value = 0.2599907647526892
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 0.197151114016889
$ws.Range("R7").Value = 1.774360026152
$ws.Range("S7").Value = 0.2437524218379954
$ws.Range("T7").Value = 0.2437524218379954

# Row 8
$ws.Range("G8").Value = 0.1935656666666667
$ws.Range("H8").Value = 0.580697
$ws.Range("I8").Value = 0.5985357597846208
$ws.Range("J8").Value = 0.5985357597846208
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.004507370114000001
$ws.Range("R8").Value = 0.040566331026
$ws.Range("S8").Value = 0.005572793168764862
$ws.Range("T8").Value = 0.005572793168764863

# Row 9
$ws.Range("G9").Value = 0.1935656666666667
$ws.Range("H9").Value = 0.580697
$ws.Range("I9").Value = 0.5985357597846208
$ws.Range("J9").Value = 0.5985357597846208
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.02572861936955555
$ws.Range("R9").Value = 0.231557574326
$ws.Range("S9").Value = 0.03181018435097395
$ws.Range("T9").Value = 0.03181018435097396

# Row 10
$ws.Range("G10").Value = 0.1935656666666667
$ws.Range("H10").Value = 0.580697
$ws.Range("I10").Value = 0.5985357597846208
$ws.Range("J10").Value = 0.5985357597846208
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 0.453869936236889
$ws.Range("R10").Value = 4.084829426132001
$ws.Range("S10").Value = 0.5611527822648821
$ws.Range("T10").Value = 0.5611527822648821

Write-Host "Updated TPM values in rows 2-10"